$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25 (pushes existing rows 25-41 down to 26-42),
# then fill it in with the new "dialog_invest_fail" entry.
$ws.Rows("25:25").Insert()

$ws.Range("A25").Value = "dialog_invest_fail"
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "name_governor"
$ws.Range("F25").Value = 0

# Append a brand-new row 43 for "dialog_not_implement_yet".
$ws.Range("A43").Value = "dialog_not_implement_yet"
$ws.Range("B43").Value = 4
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0

# Update the view to match the author's final cursor/scroll position.
try {
    $excel.ActiveWindow.ScrollRow = 13
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Viewport scroll position is cosmetic only; ignore if unsupported.
}
$ws.Range("A43").Select()
